$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "afinn"

$ws.Range("B12").Value = "NaN"

$ws.Range("C12").Value = "NaN"

$ws.Range("D12").Value = 1
$ws.Range("D12").Style = "Good"

$ws.Range("E12").Value = [double]"9.339599999999999E-44"
$ws.Range("E12").NumberFormat = "0.00E+00"

$ws.Range("F12").Value = 1
$ws.Range("F12").Style = "Good"

$ws.Range("G12").Value = [double]"1.0468E-48"
$ws.Range("G12").NumberFormat = "0.00E+00"

$ws.Range("H12").Value = 1
$ws.Range("H12").Style = "Good"

$ws.Range("I12").Value = [double]"1.5352E-76"
$ws.Range("I12").NumberFormat = "0.00E+00"

$ws.Range("J12").Value = 1
$ws.Range("J12").Style = "Good"

$ws.Range("K12").Value = [double]"3.8665000000000002E-13"
$ws.Range("K12").NumberFormat = "0.00E+00"

$ws.Range("L12").Value = [double]"0"

$ws.Range("M12").Value = [double]"0.37559999999999999"

$ws.Range("N12").Value = 1
$ws.Range("N12").Style = "Good"

$ws.Range("O12").Value = [double]"1.0890999999999999E-104"
$ws.Range("O12").NumberFormat = "0.00E+00"

$ws.Range("P12").Value = 1
$ws.Range("P12").Style = "Good"

$ws.Range("Q12").Value = [double]"6.7100000000000001E-60"
$ws.Range("Q12").NumberFormat = "0.00E+00"

$ws.Range("R12").Value = 1
$ws.Range("R12").Style = "Good"

$ws.Range("S12").Value = [double]"2.8983000000000002E-28"
$ws.Range("S12").NumberFormat = "0.00E+00"

$ws.Range("T12").Value = 1
$ws.Range("T12").Style = "Good"

$ws.Range("U12").Value = [double]"3.3134999999999998E-38"
$ws.Range("U12").NumberFormat = "0.00E+00"

$ws.Range("V12").Value = 1
$ws.Range("V12").Style = "Good"

$ws.Range("W12").Value = [double]"2.9907999999999999E-5"
$ws.Range("W12").NumberFormat = "0.00E+00"

$ws.Range("X12").Value = 1
$ws.Range("X12").Style = "Good"

$ws.Range("Y12").Value = [double]"4.4576000000000002E-98"
$ws.Range("Y12").NumberFormat = "0.00E+00"

$ws.Range("Z12").Value = [double]"0"

$ws.Range("AA12").Value = [double]"7.7499999999999999E-2"

$ws.Range("AB12").Value = 1
$ws.Range("AB12").Style = "Good"

$ws.Range("AC12").Value = [double]"3.9390999999999997E-4"
$ws.Range("AC12").NumberFormat = "0.00E+00"

$ws.Range("AD12").Value = [double]"0"

$ws.Range("AE12").Value = [double]"0.2437"

$ws.Range("H16").Select()

Write-Host "done"